# Update crypto price/volume data per Sat Mar 30 07:07:34 UTC 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 48/49 swap coins (B,C columns) in addition to the usual D/E numeric updates.
$updates = @(
    @{ Cell = 'D2'; Value = '69.873.97' },
    @{ Cell = 'E2'; Value = '  -0.80%  ' },
    @{ Cell = 'D3'; Value = '3.499.91' },
    @{ Cell = 'E3'; Value = '  -2.05%  ' },
    @{ Cell = 'E4'; Value = '  -0.14%  ' },
    @{ Cell = 'D5'; Value = '603.14' },
    @{ Cell = 'E5'; Value = '  -1.10%  ' },
    @{ Cell = 'D6'; Value = '198.27' },
    @{ Cell = 'E6'; Value = '  +6.21%  ' },
    @{ Cell = 'E7'; Value = '  +0.99%  ' },
    @{ Cell = 'E8'; Value = '  -0.08%  ' },
    @{ Cell = 'E9'; Value = '  -2.94%  ' },
    @{ Cell = 'E10'; Value = '  +1.21%  ' },
    @{ Cell = 'D11'; Value = '54.41' },
    @{ Cell = 'E11'; Value = '  +0.71%  ' },
    @{ Cell = 'E12'; Value = '  -2.57%  ' },
    @{ Cell = 'E13'; Value = '  +0.72%  ' },
    @{ Cell = 'D14'; Value = '4.057.29' },
    @{ Cell = 'E14'; Value = '  -1.90%  ' },
    @{ Cell = 'D15'; Value = '593.52' },
    @{ Cell = 'E15'; Value = '  +2.07%  ' },
    @{ Cell = 'D16'; Value = '69.942.74' },
    @{ Cell = 'E16'; Value = '  -0.76%  ' },
    @{ Cell = 'D17'; Value = '19.04' },
    @{ Cell = 'E17'; Value = '  -0.06%  ' },
    @{ Cell = 'E18'; Value = '  -0.94%  ' },
    @{ Cell = 'D19'; Value = '3.508.23' },
    @{ Cell = 'E19'; Value = '  -1.35%  ' },
    @{ Cell = 'E20'; Value = '  +0.09%  ' },
    @{ Cell = 'D21'; Value = '0.988' },
    @{ Cell = 'E21'; Value = '  -0.88%  ' },
    @{ Cell = 'D22'; Value = '18.17' },
    @{ Cell = 'E22'; Value = '  +4.49%  ' },
    @{ Cell = 'D23'; Value = '104.54' },
    @{ Cell = 'E23'; Value = '  +10.38%  ' },
    @{ Cell = 'E24'; Value = '  -2.81%  ' },
    @{ Cell = 'E25'; Value = '  +2.84%  ' },
    @{ Cell = 'D26'; Value = '3.10' },
    @{ Cell = 'E26'; Value = '  +4.44%  ' },
    @{ Cell = 'D27'; Value = '10.96' },
    @{ Cell = 'E27'; Value = '  +0.08%  ' },
    @{ Cell = 'D28'; Value = '9.76' },
    @{ Cell = 'E28'; Value = '  +3.56%  ' },
    @{ Cell = 'D29'; Value = '33.58' },
    @{ Cell = 'E29'; Value = '  +3.81%  ' },
    @{ Cell = 'D30'; Value = '4.55' },
    @{ Cell = 'E30'; Value = '  +22.60%  ' },
    @{ Cell = 'D31'; Value = '7.26' },
    @{ Cell = 'E31'; Value = '  +2.54%  ' },
    @{ Cell = 'D32'; Value = '12.74' },
    @{ Cell = 'E32'; Value = '  +3.67%  ' },
    @{ Cell = 'E33'; Value = '  +0.73%  ' },
    @{ Cell = 'D34'; Value = '63.66' },
    @{ Cell = 'E34'; Value = '  +0.27%  ' },
    @{ Cell = 'D35'; Value = '3.727.86' },
    @{ Cell = 'E35'; Value = '  +5.75%  ' },
    @{ Cell = 'D36'; Value = '0.0₃0811' },
    @{ Cell = 'E36'; Value = '  +2.79%  ' },
    @{ Cell = 'E37'; Value = '  +0.01%  ' },
    @{ Cell = 'D38'; Value = '511.23' },
    @{ Cell = 'E38'; Value = '  -4.17%  ' },
    @{ Cell = 'D39'; Value = '0.391' },
    @{ Cell = 'E39'; Value = '  -3.58%  ' },
    @{ Cell = 'D40'; Value = '2.98' },
    @{ Cell = 'E40'; Value = '  -7.66%  ' },
    @{ Cell = 'D41'; Value = '36.69' },
    @{ Cell = 'E41'; Value = '  -1.84%  ' },
    @{ Cell = 'D42'; Value = '3.53' },
    @{ Cell = 'E42'; Value = '  -0.01%  ' },
    @{ Cell = 'E43'; Value = '  +0.43%  ' },
    @{ Cell = 'D44'; Value = '0.0458' },
    @{ Cell = 'E44'; Value = '  -0.83%  ' },
    @{ Cell = 'E45'; Value = '  -3.43%  ' },
    @{ Cell = 'E46'; Value = '  -1.10%  ' },
    @{ Cell = 'D47'; Value = '3.32' },
    @{ Cell = 'E47'; Value = '  -2.76%  ' },
    @{ Cell = 'B48'; Value = 'THORChain' },
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune' },
    @{ Cell = 'D48'; Value = '8.76' },
    @{ Cell = 'E48'; Value = '  -5.30%  ' },
    @{ Cell = 'B49'; Value = 'FirstDigitalUSD' },
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd' },
    @{ Cell = 'D49'; Value = '1.00' },
    @{ Cell = 'E49'; Value = '  +0.36%  ' },
    @{ Cell = 'D50'; Value = '132.42' },
    @{ Cell = 'E50'; Value = '  -3.01%  ' },
    @{ Cell = 'D51'; Value = '0.000242' },
    @{ Cell = 'E51'; Value = '  -2.21%  ' }
)

foreach ($u in $updates) {
    $cell = $u.Cell
    $value = $u.Value
    $range = $ws.Range($cell)

    if ($cell -match '^D') {
        # Column D holds price strings that look numeric (e.g. "3.10", "0.988",
        # "69.873.97"). Excel's automatic type coercion would turn these into
        # real numbers (dropping formatting like trailing zeros), so force a
        # text entry via the classic leading-apostrophe trick, then strip the
        # quote-prefix cell style back off so no stray formatting is left
        # behind on the cell.
        $range.Value = "'" + $value
        $range.Style = "Normal"
    } else {
        # Columns B, C and E (coin name, link, percent change) are safe to
        # assign directly - none of them are ambiguous with Excel's numeric
        # auto-detection.
        $range.Value = $value
    }
}

Write-Output "Updated $($updates.Count) cells"
